$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New TikTok search results shift into rows 4-7 (row 4's old data scrolled off,
# row 5 moved up to row 4, and three fresh results were appended).
# Video IDs are large integers that must stay exact text (not lossy floats),
# so force Text number format before writing them.

$ws.Range("C4:C7").NumberFormat = "@"

$ws.Range("A4").Value = "https://www.tiktok.com/@wildmane_mainecoons/video/7544160476519779615"
$ws.Range("B4").Value = "wildmane_mainecoons"
$ws.Range("C4").Value = "7544160476519779615"
$ws.Range("D4").Value = "Video by @wildmane_mainecoons"

$ws.Range("A5").Value = "https://www.tiktok.com/@cute.catsxxx/video/7543767375758265655"
$ws.Range("B5").Value = "cute.catsxxx"
$ws.Range("C5").Value = "7543767375758265655"
$ws.Range("D5").Value = "Video by @cute.catsxxx"

$ws.Range("A6").Value = "https://www.tiktok.com/@catutucom/video/7543337508155804935"
$ws.Range("B6").Value = "catutucom"
$ws.Range("C6").Value = "7543337508155804935"
$ws.Range("D6").Value = "Video by @catutucom"

$ws.Range("A7").Value = "https://www.tiktok.com/@icecreamtina0/video/7543814939123666196"
$ws.Range("B7").Value = "icecreamtina0"
$ws.Range("C7").Value = "7543814939123666196"
$ws.Range("D7").Value = "Video by @icecreamtina0"
